# fitur : hapus data parkir selesai
# Updates the "Data Parkir" row for vehicle 6CFX94LN9BM1 with its actual
# checkout time/duration/payment, removes the now-redundant "still parked"
# row (PCARMF5WXU87), and refreshes the "Ringkasan Keuangan" summary sheet
# to reflect the corrected totals.

$wb = $excel.ActiveWorkbook

# --- Sheet: Data Parkir ---
$ws1 = $wb.Worksheets.Item("Data Parkir")

# Row 4 (6CFX94LN9BM1): update exit time, duration and amount paid
$ws1.Range("E4").Value = "2025-02-03 17:44:07"
$ws1.Range("F4").Value = "00:04:01"
$ws1.Range("H4").Value = 50000

# Row 5 (PCARMF5WXU87, still parked) is no longer needed - remove it entirely
$ws1.Rows(5).Delete()

# --- Sheet: Ringkasan Keuangan ---
$ws2 = $wb.Worksheets.Item("Ringkasan Keuangan")

$ws2.Range("B2").Value = "Rp 8,000"
$ws2.Range("B3").Value = "Rp 107,000"
$ws2.Range("B4").Value = "Rp 8,000"
$ws2.Range("B5").Value = "Delapan ribu Rupiah"
